# Delete duplicate fastq file entry: remove entire row 2 (the first data row),
# shifting all subsequent rows up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Select() | Out-Null
$ws.Rows.Item(2).Delete()
